$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Test Description for row 14 (Validate_validPinCodeAtRegistrationPage)
# with the fuller, multi-point description. Setting this value causes the old
# shared string to be replaced / reordered so that "Share" (used by row 13,
# Application_Share_Through_WhatsApp) and the new description end up exactly
# where the target workbook expects them.
$newDescription = "verify & validate that application is checking the below mentioned points;`n- valid and invalid pin code format in registration page.`n- length allowed for zipcode field.`n- Fetching dynamic zipcodes from Postal APi.`n- Parsing one by one into address fields."
$ws.Range("D14").Value = $newDescription

# The longer, multi-line text needs a taller row to display properly.
$ws.Rows.Item(14).RowHeight = 90

# Update the view: scroll so row 6 is at the top and select C14.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
